$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (values changed per the revision) ---
# Row 2
$ws.Range("D2").Value = 44904
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 1500
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44904
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 1600
$ws.Range("O4").Value = 1600
$ws.Range("P4").Value = 1600
$ws.Range("Q4").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S4").Value = 1600
$ws.Range("T4").Value = 1

# Row 6
$ws.Range("D6").Value = 44343
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 1700
$ws.Range("O6").Value = 1700
$ws.Range("P6").Value = 1700
$ws.Range("S6").Value = 1700

# Row 7
$ws.Range("D7").Value = 44336
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 1500
$ws.Range("O7").Value = 1500
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44371
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 1800
$ws.Range("O8").Value = 1800
$ws.Range("P8").Value = 1800
$ws.Range("Q8").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S8").Value = 1800
$ws.Range("T8").Value = 1

# Row 9
$ws.Range("D9").Value = 44371
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 1200
$ws.Range("O9").Value = 1200
$ws.Range("P9").Value = 1200
$ws.Range("S9").Value = 1200

# Row 10
$ws.Range("D10").Value = 45113
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 6
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("Q10").Value = '$/bandeja 10 kilos'
$ws.Range("S10").Value = 2400
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 45113
$ws.Range("M11").Value = 8

# Row 12
$ws.Range("D12").Value = 45113
$ws.Range("M12").Value = 15

# Row 13
$ws.Range("D13").Value = 45113
$ws.Range("L13").Value = 'Tercera'
$ws.Range("M13").Value = 8
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = '$/bandeja 10 kilos'
$ws.Range("S13").Value = 1200
$ws.Range("T13").Value = 10

# Row 14
$ws.Range("D14").Value = 44880
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value = 2000
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44880
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 180
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("S15").Value = 1500

# --- Append new rows 16-19 ---
# Row 16
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C16").Value = 'Metropolitana'
$ws.Range("D16").Value = 44195
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = 'Tropicales y subtropicales'
$ws.Range("I16").Value = 100108004
$ws.Range("J16").Value = 'Papaya'
$ws.Range("K16").Value = 'Cultivar IV Región'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = '$/bandeja 10 kilos'
$ws.Range("R16").Value = 'Provincia del Elquí'
$ws.Range("S16").Value = 1500
$ws.Range("T16").Value = 10

# Row 17
$ws.Range("A17").Value = 9
$ws.Range("B17").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C17").Value = 'Metropolitana'
$ws.Range("D17").Value = 44292
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = 'Tropicales y subtropicales'
$ws.Range("I17").Value = 100108004
$ws.Range("J17").Value = 'Papaya'
$ws.Range("K17").Value = 'Cultivar IV Región'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 14000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = 'Provincia del Elquí'
$ws.Range("S17").Value = 1400
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C18").Value = 'Metropolitana'
$ws.Range("D18").Value = 44391
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 'Fruta'
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = 'Tropicales y subtropicales'
$ws.Range("I18").Value = 100108004
$ws.Range("J18").Value = 'Papaya'
$ws.Range("K18").Value = 'Cultivar IV Región'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 15
$ws.Range("N18").Value = 1500
$ws.Range("O18").Value = 1500
$ws.Range("P18").Value = 1500
$ws.Range("Q18").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R18").Value = 'Provincia del Elquí'
$ws.Range("S18").Value = 1500
$ws.Range("T18").Value = 1

# Row 19
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C19").Value = 'Metropolitana'
$ws.Range("D19").Value = 44391
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 'Fruta'
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = 'Tropicales y subtropicales'
$ws.Range("I19").Value = 100108004
$ws.Range("J19").Value = 'Papaya'
$ws.Range("K19").Value = 'Cultivar IV Región'
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 20
$ws.Range("N19").Value = 1000
$ws.Range("O19").Value = 1000
$ws.Range("P19").Value = 1000
$ws.Range("Q19").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R19").Value = 'Provincia del Elquí'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 1

# Apply the date/time number format to the newly added "Fecha" cells so they match
# the existing date column formatting (style used by D2:D15).
$ws.Range("D16:D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
